# Update the cryptos.xlsx price/volume snapshot (GitHub Actions refresh).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). D/E are stored as plain
# text (not numbers), matching the original inlineStr cells. For D-column
# values that look like a bare number (e.g. "0.999", "93.05") we force the
# cell to Text format first so Excel's COM layer doesn't silently coerce
# the assignment into a numeric cell - this mirrors the real "type text
# into a Text-formatted cell" behaviour and keeps the stored value a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.724.53"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.472.08"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.58"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.05"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0876"
$ws.Range("E10").Value = "  +10.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.06"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D13").Value = "2.855.33"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.58"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "2.446.52"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.802"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "41.684.97"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.14"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.09"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.07"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.75"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.06"
$ws.Range("E30").Value = "  +4.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.15"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0764"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.34"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  +7.71%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.02"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "1.996.77"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.45"
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("D48").Value = "2.716.96"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.13"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.38"
$ws.Range("E50").Value = "  +6.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.21"
$ws.Range("E51").Value = "  +0.00%  "
